$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the three data rows before writing the
# updated exposure-site records back in.
$ws.Range("A2:E4").ClearContents()

# Row 2: replace the Black Rock exposure-site entry with the McKinnon one
$ws.Range("A2").Value = "McKinnon"
$ws.Range("B2").Value = "Hotlocks By Rachael Hairdresser, 260 McKinnon Road, McKinnon VIC 3204"
$ws.Range("C2").Value = "23/12/20 4:00pm-6:00pm"
$ws.Range("D2").Value = "Case had hair cut in store"
$ws.Range("E2").Value = "new"

# Row 3: the old "Left Bank Melbourne" entry is dropped; the Melbourne Central
# Lion Hotel entry (previously row 4) moves up and is now marked "old"
$ws.Range("A3").Value = "Melbourne"
$ws.Range("B3").Value = "Melbourne Central Lion Hotel, 211 La Trobe Street"
$ws.Range("C3").Value = "28/12/2020 10pm - 12.30am"
$ws.Range("D3").Value = "Case attended Venue"
$ws.Range("E3").Value = "old"

# Row 4: a new Southbank restaurant entry is added
$ws.Range("A4").Value = "Southbank"
$ws.Range("B4").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C4").Value = "23/12/20 8:00pm-11:00pm"
$ws.Range("D4").Value = "Case attended restaurant"
$ws.Range("E4").Value = "new"

# Widen columns B and D to fit the longer replacement text, and match the
# final selection recorded in the saved file
$ws.Columns.Item(2).ColumnWidth = 59
$ws.Columns.Item(4).ColumnWidth = 19.833333333333332
$ws.Range("C4").Select() | Out-Null
